$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "09-09-2021"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 3500
$ws.Range("C6").Value = 5870
$ws.Range("D6").Value = 3500
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 3420
$ws.Range("G6").Value = 2.55
